$d = $word.ActiveDocument

# Remove every paragraph except the final (already-empty) one that sits
# right before the sectPr mark. The source-link content (the
# "Adatbázisok:" heading, the "Kaggle.com:" sub-heading, and the five
# dataset reference lines) was moved out to its own file, so the body
# collapses down to a single empty paragraph.
while ($d.Paragraphs.Count -gt 1) {
    $d.Content.Delete()
}
